$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: push a pure-numeric-looking string (e.g. "37%") into a target cell as
# literal TEXT. A direct $ws.Range(cell).Value = "37%" gets auto-parsed by Excel as
# the number 0.37 with a new percentage number format (changing both the stored
# type and the cell's style). Instead, build the text on an out-of-range scratch
# cell via a ="..." formula (so no quotePrefix/style is minted either), copy it,
# and PasteSpecial only the *value* (xlPasteValues = -4163) onto the target so the
# target keeps its original style untouched. Then clear the scratch cell.
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range("E2").Value = "2026-02-23 22:48:26"
Set-TextValue "H2" "37%"
$ws.Range("O2").Value = "5.6 °C"
$ws.Range("E3").Value = "2026-02-23 22:48:28"
$ws.Range("E4").Value = "2026-02-23 22:48:30"
$ws.Range("J4").Value = "1024.6 hPa"
$ws.Range("E5").Value = "2026-02-23 22:48:32"
Set-TextValue "H5" "28%"
$ws.Range("E6").Value = "2026-02-23 22:48:35"
$ws.Range("E7").Value = "2026-02-23 22:48:37"
$ws.Range("E8").Value = "2026-02-23 22:48:39"
$ws.Range("E9").Value = "2026-02-23 22:48:41"
Set-TextValue "H9" "74%"
$ws.Range("O9").Value = "12.3 °C"
$ws.Range("E10").Value = "2026-02-23 22:48:44"
$ws.Range("O10").Value = "10.5 °C"
$ws.Range("E11").Value = "2026-02-23 22:48:46"
Set-TextValue "H11" "70%"
$ws.Range("O11").Value = "8.6 °C"
$ws.Range("E12").Value = "2026-02-23 22:48:48"
$ws.Range("E13").Value = "2026-02-23 22:48:51"
Set-TextValue "H13" "60%"
$ws.Range("K13").Value = "15.6 MJ/m2"
$ws.Range("O13").Value = "6.9 °C"
$ws.Range("E14").Value = "2026-02-23 22:48:54"
$ws.Range("O14").Value = "12.2 °C"
$ws.Range("E15").Value = "2026-02-23 22:48:55"
$ws.Range("O15").Value = "12.2 °C"
$ws.Range("E16").Value = "2026-02-23 22:48:56"
$ws.Range("E17").Value = "2026-02-23 22:48:57"
Set-TextValue "H17" "42%"
$ws.Range("O17").Value = "8.5 °C"
$ws.Range("E18").Value = "2026-02-23 22:48:58"
Set-TextValue "H18" "75%"
$ws.Range("O18").Value = "10.7 °C"
$ws.Range("E19").Value = "2026-02-23 22:48:59"
$ws.Range("O19").Value = "12.2 °C"
$ws.Range("E20").Value = "2026-02-23 22:49:02"
$ws.Range("E21").Value = "2026-02-23 22:49:04"
Set-TextValue "H21" "60%"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-23 22:49:06"
$ws.Range("E23").Value = "2026-02-23 22:49:09"
Set-TextValue "H23" "23%"
$ws.Range("E24").Value = "2026-02-23 22:49:11"
$ws.Range("E25").Value = "2026-02-23 22:49:13"
$ws.Range("E26").Value = "2026-02-23 22:49:16"
$ws.Range("E27").Value = "2026-02-23 22:49:19"
$ws.Range("K27").Value = "16.8 MJ/m2"
$ws.Range("E28").Value = "2026-02-23 22:49:21"
$ws.Range("E29").Value = "2026-02-23 22:49:23"
Set-TextValue "H29" "84%"
$ws.Range("O29").Value = "10.5 °C"
$ws.Range("E30").Value = "2026-02-23 22:49:25"
$ws.Range("E31").Value = "2026-02-23 22:49:28"
$ws.Range("E32").Value = "2026-02-23 22:49:31"
Set-TextValue "H32" "69%"
$ws.Range("N32").Value = "-1.5 °C 22:14 TU"
$ws.Range("O32").Value = "7.0 °C"
$ws.Range("E33").Value = "2026-02-23 22:49:33"
$ws.Range("E34").Value = "2026-02-23 22:49:36"
Set-TextValue "H34" "44%"
$ws.Range("O34").Value = "3.9 °C"
$ws.Range("E35").Value = "2026-02-23 22:49:38"
$ws.Range("O35").Value = "11.9 °C"
$ws.Range("E36").Value = "2026-02-23 22:49:41"
$ws.Range("K36").Value = "15.2 MJ/m2"
$ws.Range("E37").Value = "2026-02-23 22:49:44"
$ws.Range("J37").Value = "1026.8 hPa"
$ws.Range("O37").Value = "8.8 °C"
$ws.Range("E38").Value = "2026-02-23 22:49:46"
$ws.Range("K38").Value = "15.8 MJ/m2"
$ws.Range("E39").Value = "2026-02-23 22:49:48"
Set-TextValue "H39" "27%"
$ws.Range("E40").Value = "2026-02-23 22:49:51"
$ws.Range("O40").Value = "8.5 °C"
$ws.Range("E41").Value = "2026-02-23 22:49:53"
$ws.Range("J41").Value = "1024.6 hPa"
$ws.Range("O41").Value = "11.7 °C"
$ws.Range("E42").Value = "2026-02-23 22:49:56"
$ws.Range("E43").Value = "2026-02-23 22:49:58"
Set-TextValue "H43" "71%"
$ws.Range("E44").Value = "2026-02-23 22:50:01"
$ws.Range("O44").Value = "3.0 °C"
$ws.Range("E45").Value = "2026-02-23 22:50:03"
$ws.Range("O45").Value = "8.1 °C"
$ws.Range("E46").Value = "2026-02-23 22:50:06"
Set-TextValue "H46" "74%"
$ws.Range("O46").Value = "10.0 °C"
